# Update on 20250607 part 2
# Insert a new "山东" (Shandong) script entry (shandongsd.js) as the first
# row of the existing 山东 group on sheet "地方台JS脚本" (sheet 1), pushing
# the previous single-row 山东 entry (jinan.js) down to become the second
# (merged-continuation) row of a new two-row 山东 group.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Insert a new row above the current row 40 ("山东" / jinan.js). This
#    shifts rows 40-62 down to 41-63 and keeps all existing cell content,
#    styles and merged cells intact (Excel re-points them automatically).
$ws.Rows.Item(40).Insert()

# 2. Populate the freshly inserted row 40 with the new Shandong entry.
$ws.Range("A40").Value = "山东"
$ws.Range("B40").Value = "shandongsd.js"
$ws.Range("C40").Value = "本地"
$ws.Range("D40").Value = "卫视,山东省频道"

# 3. Copy cell formatting from already-correctly-styled neighbour cells so
#    we reuse existing style/border definitions instead of inventing new
#    ones. A40 should look like the top cell of a merged "province" block
#    (border on left/right/top only); row 42 column A (old row 41, the
#    "山西" header) already has exactly that style.
$ws.Range("A42").Copy()
$ws.Range("A40").PasteSpecial(-4122)

# B/C/D on row 40 should use the plain bordered-box style already used by
# B41/C41/D41 (old row 40 cells, now shifted down to row 41).
$ws.Range("B41").Copy()
$ws.Range("B40").PasteSpecial(-4122)
$ws.Range("C41").Copy()
$ws.Range("C40").PasteSpecial(-4122)
$ws.Range("D41").Copy()
$ws.Range("D40").PasteSpecial(-4122)

# 4. Row 41 (previously row 40, "jinan.js") becomes the bottom/continuation
#    row of the merged "山东" block: clear its A-column text and restyle it
#    like the bottom cell of a merged province block (e.g. the existing
#    bottom-of-merge cell at A49, old A48 "汾阳市频道" header of 山西 group).
$ws.Range("A41").ClearContents()
$ws.Range("A49").Copy()
$ws.Range("A41").PasteSpecial(-4122)

# 5. Merge A40:A41 into a single "省份/直辖市" cell, matching the rest of
#    the sheet's merged province column.
$ws.Range("A40:A41").Merge()

$app.CutCopyMode = 0

# 6. Fix up ranges that Insert() does not automatically extend: the
#    worksheet AutoFilter and the hidden _FilterDatabase defined name both
#    still reference the old last row (62) and need to cover row 63 now.
#    Range.AutoFilter() toggles filtering off when it is already on, so we
#    must explicitly drop it first before re-applying over the new range.
$ws.AutoFilterMode = $false
$ws.Range("A1:D63").AutoFilter()

$ws.Names.Item(1).RefersTo = "=地方台JS脚本!`$A`$1:`$D`$63"
